$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price column cells we touch so Excel COM
# does not reinterpret numeric-looking strings (e.g. "1.005") as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.037.31'
$ws.Range("E2").Value = '  -2.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.798.66'
$ws.Range("E3").Value = '  -2.49%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.10'
$ws.Range("E5").Value = '  -2.42%  '

$ws.Range("E6").Value = '  +0.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4218'
$ws.Range("E7").Value = '  -2.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3602'
$ws.Range("E8").Value = '  -2.81%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07261'
$ws.Range("E9").Value = '  -1.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8426'
$ws.Range("E10").Value = '  -3.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.29'
$ws.Range("E11").Value = '  -3.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.824.79'
$ws.Range("E12").Value = '  -5.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.298'
$ws.Range("E13").Value = '  -3.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.376'
$ws.Range("E14").Value = '  -3.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06779'
$ws.Range("E15").Value = '  -2.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.008'
$ws.Range("E16").Value = '  +0.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.69'
$ws.Range("E17").Value = '  -0.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008760'
$ws.Range("E18").Value = '  -3.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.06'
$ws.Range("E20").Value = '  -3.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.211.44'
$ws.Range("E21").Value = '  -2.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.084'
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.07'
$ws.Range("E23").Value = '  +0.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.064.48'
$ws.Range("E24").Value = '  -4.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.928'
$ws.Range("E25").Value = '  -2.98%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.19'
$ws.Range("E26").Value = '  -0.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.15'
$ws.Range("E27").Value = '  -4.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.019'
$ws.Range("E28").Value = '  -6.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.58'
$ws.Range("E29").Value = '  -1.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.656'
$ws.Range("E30").Value = '  -12.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09008'
$ws.Range("E31").Value = '  +1.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7302'
$ws.Range("E32").Value = '  -7.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.866'
$ws.Range("E33").Value = '  -3.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.347'
$ws.Range("E34").Value = '  -5.92%  '

$ws.Range("E35").Value = '  -6.73%  '

$ws.Range("E36").Value = '  +0.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.081'
$ws.Range("E37").Value = '  -2.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05154'
$ws.Range("E38").Value = '  -5.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01905'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4999'
$ws.Range("E40").Value = '  -3.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1633'
$ws.Range("E41").Value = '  -3.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.627'
$ws.Range("E42").Value = '  -7.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.092'
$ws.Range("E43").Value = '  -6.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.962'
$ws.Range("E44").Value = '  -12.41%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.21'
$ws.Range("E45").Value = '  -1.48%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.27'
$ws.Range("E46").Value = '  -3.83%  '

$ws.Range("E47").Value = '  +0.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06313'
$ws.Range("E48").Value = '  -3.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4530'
$ws.Range("E49").Value = '  -5.37%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.602'
$ws.Range("E50").Value = '  -3.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.719'
$ws.Range("E51").Value = '  -8.00%  '
